# Add "Wins", "Losses", "Ties" (season record) columns to the worksheet.
# This mirrors the commit's addition of AD/AE/AF columns: headers in row 1
# (styled like the rest of the header row) and the team's season record
# (64 wins, 98 losses, 0 ties) repeated for every player row (2-46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the new
# header cells so AD1:AF1 pick up the same bold/centered/bordered style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$wins = 64
$losses = 98
$ties = 0

for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}
